$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row for Mayo (May) readings
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Mayo"
$ws.Range("C9").Value = 2024
$ws.Range("D9").Value = 1042
$ws.Range("E9").Value = 120

# Match the styles used by the other data rows (A,C,D,E use style index 1 / "s=1",
# B uses style index 2 / "s=2"); copy formatting from row 8 which already has the
# correct number formats applied.
$ws.Range("A8:E8").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122) # xlPasteFormats

# Update the active selection to mirror the saved workbook view state
$ws.Range("G12").Select()
